# Applies the changes described by the commit diff:
#  - Slide 2 ("Основные Возможности"): wrap "зелёных задач" and
#    "красных дистракторов" in curly quotes.
#  - Slide 5 ("Класс Player ..."): drop the literal key hint
#    "(W, A, S, D)" from the movement bullet.
#  - Slide 7 (last slide): remove the leftover screenshot/picture
#    shape at the bottom of the slide.

$p = $ppt.ActivePresentation

# --- Slide 2: quote the "зелёных задач" / "красных дистракторов" bullets ---
$slide2 = $p.Slides.Item(2)
$body2 = $slide2.Shapes.Item(2)
$tr2 = $body2.TextFrame.TextRange

$para2a = $tr2.Paragraphs(2, 1)
$run2a = $para2a.Runs(1, 1)
$run2a.Text = "Сбор “зелёных задач” (+10 очков)"

$para2b = $tr2.Paragraphs(3, 1)
$run2b = $para2b.Runs(1, 1)
$run2b.Text = "Уничтожение “красных дистракторов” (+8 очков)"

# --- Slide 5: trim the "(W, A, S, D)" suffix from the movement bullet ---
$slide5 = $p.Slides.Item(5)
$body5 = $slide5.Shapes.Item(2)
$tr5 = $body5.TextFrame.TextRange

$para5 = $tr5.Paragraphs(2, 1)
$run5 = $para5.Runs(1, 1)
$run5.Text = "Движение во всех направлениях"

# --- Slide 7: delete the stray picture shape at the bottom of the slide ---
$slide7 = $p.Slides.Item(7)
$slide7.Shapes.Item(3).Delete()
